$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Aggiornamento dati fino al 13/05 (righe 252-255)
# Replica lo stile (data, grassetto, bordo, centrato) della colonna A
# dall'ultima riga esistente (251) alle nuove righe.
$ws.Range("A251").Copy($ws.Range("A252:A255"))

$ws.Range("A252").Value = 44326
$ws.Range("B252").Value = 0
$ws.Range("C252").Value = 16
$ws.Range("D252").Value = 66.32949175026947

$ws.Range("A253").Value = 44327
$ws.Range("B253").Value = 1
$ws.Range("C253").Value = 15
$ws.Range("D253").Value = 62.18389851587763

$ws.Range("A254").Value = 44328
$ws.Range("B254").Value = 1
$ws.Range("C254").Value = 15
$ws.Range("D254").Value = 62.18389851587763

$ws.Range("A255").Value = 44329
$ws.Range("B255").Value = 2
$ws.Range("C255").Value = 16
$ws.Range("D255").Value = 66.32949175026947
